$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''27.901.80'
$ws.Range("E2").Value = '  -2.46%  '

$ws.Range("D3").Value = '''1.794.76'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("D4").Value = '''0.9994'
$ws.Range("E4").Value = '  -0.27%  '

$ws.Range("D5").Value = '''316.68'
$ws.Range("E5").Value = '  +0.00%  '

$ws.Range("D6").Value = '''0.9994'
$ws.Range("E6").Value = '  -0.22%  '

$ws.Range("D7").Value = '''0.5326'

$ws.Range("D8").Value = '''0.3877'
$ws.Range("E8").Value = '  +3.13%  '

$ws.Range("D9").Value = '''0.07454'
$ws.Range("E9").Value = '  -0.68%  '

$ws.Range("D10").Value = '''41.37'
$ws.Range("E10").Value = '  -2.46%  '

$ws.Range("D11").Value = '''1.087'
$ws.Range("E11").Value = '  -2.45%  '

$ws.Range("D12").Value = '''0.9989'
$ws.Range("E12").Value = '  -0.34%  '

$ws.Range("D13").Value = '''6.190'
$ws.Range("E13").Value = '  +0.65%  '

$ws.Range("D14").Value = '''7.437'
$ws.Range("E14").Value = '  +0.78%  '

$ws.Range("D15").Value = '''20.34'
$ws.Range("E15").Value = '  -1.71%  '

$ws.Range("D16").Value = '''1.790.40'
$ws.Range("E16").Value = '  -0.52%  '

$ws.Range("D17").Value = '''88.37'
$ws.Range("E17").Value = '  -2.04%  '

$ws.Range("E18").Value = '  -0.55%  '

$ws.Range("D19").Value = '''0.06538'
$ws.Range("E19").Value = '  +1.32%  '

$ws.Range("D20").Value = '''0.9996'
$ws.Range("E20").Value = '  -0.12%  '

$ws.Range("D21").Value = '''17.25'
$ws.Range("E21").Value = '  +0.04%  '

$ws.Range("E22").Value = '  +0.82%  '

$ws.Range("D23").Value = '''27.927.28'
$ws.Range("E23").Value = '  -2.45%  '

$ws.Range("E24").Value = '  +0.03%  '

$ws.Range("E25").Value = '  -0.21%  '

$ws.Range("D26").Value = '''156.37'
$ws.Range("E26").Value = '  -1.31%  '

$ws.Range("D27").Value = '''20.15'
$ws.Range("E27").Value = '  -1.41%  '

$ws.Range("D28").Value = '''1.998.20'
$ws.Range("E28").Value = '  -0.42%  '

$ws.Range("D29").Value = '''2.301'

$ws.Range("D30").Value = '''121.76'
$ws.Range("E30").Value = '  -1.03%  '

$ws.Range("D31").Value = '''0.1090'
$ws.Range("E31").Value = '  +3.15%  '

$ws.Range("D32").Value = '''1.099'
$ws.Range("E32").Value = '  -0.42%  '

$ws.Range("E33").Value = '  -0.48%  '

$ws.Range("D34").Value = '''5.507'
$ws.Range("E34").Value = '  -2.45%  '

$ws.Range("D35").Value = '''0.06975'
$ws.Range("E35").Value = '  +7.33%  '

$ws.Range("D36").Value = '''0.2204'
$ws.Range("E36").Value = '  -1.74%  '

$ws.Range("D37").Value = '''0.02272'
$ws.Range("E37").Value = '  -1.26%  '

$ws.Range("D38").Value = '''5.070'
$ws.Range("E38").Value = '  +0.53%  '

$ws.Range("D39").Value = '''8.400'
$ws.Range("E39").Value = '  -3.70%  '

$ws.Range("D40").Value = '''11.25'
$ws.Range("E40").Value = '  +0.04%  '

$ws.Range("E41").Value = '  -0.55%  '

$ws.Range("D42").Value = '''0.6119'
$ws.Range("E42").Value = '  -1.63%  '

$ws.Range("D43").Value = '''1.413'
$ws.Range("E43").Value = '  -0.72%  '

$ws.Range("D44").Value = '''13.25'
$ws.Range("E44").Value = '  +0.12%  '

$ws.Range("D45").Value = '''3.675'
$ws.Range("E45").Value = '  -0.36%  '

$ws.Range("D46").Value = '''0.5713'
$ws.Range("E46").Value = '  -2.30%  '

$ws.Range("D47").Value = '''124.65'
$ws.Range("E47").Value = '  -1.18%  '

$ws.Range("D48").Value = '''1.913'
$ws.Range("E48").Value = '  -1.41%  '

$ws.Range("D49").Value = '''1.176'
$ws.Range("E49").Value = '  +1.89%  '

$ws.Range("D50").Value = '''0.06797'
$ws.Range("E50").Value = '  -1.37%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").Value = '''0.00000000293'
$ws.Range("E51").Value = '  +32.62%  '
